$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.692005395889282
$ws.Range("E2").Value = 5865.199779030891
$ws.Range("F2").Value = 0.1975374015700881
$ws.Range("G2").Value = 0.1809931001826227
$ws.Range("H2").Value = 0.1695156604601038
$ws.Range("I2").Value = 0.1446706062871062
$ws.Range("J2").Value = 0.1335253691896111
$ws.Range("K2").Value = 0.1280737292330381
$ws.Range("L2").Value = 0.1240173867295842
$ws.Range("M2").Value = 0.1224014423797636
$ws.Range("N2").Value = 0.1223416702252344
$ws.Range("O2").Value = 0.1223416702252344
$ws.Range("P2").Value = 0.1223416702252344
$ws.Range("Q2").Value = 0.1223416702252344
$ws.Range("R2").Value = 0.1223313797082045
$ws.Range("S2").Value = 0.1223313797082045
$ws.Range("T2").Value = 0.1223313797082045
$ws.Range("U2").Value = 0.1223313797082045
$ws.Range("V2").Value = 0.1223313797082045
$ws.Range("W2").Value = 0.1223313797082045
$ws.Range("X2").Value = 0.1223313797082045
$ws.Range("Y2").Value = 0.1223313797082045
$ws.Range("C3").Value = 1.660001277923584
$ws.Range("E3").Value = 5864.944556532597
$ws.Range("F3").Value = 0.2029301912968721
$ws.Range("G3").Value = 0.1793876553257469
$ws.Range("H3").Value = 0.1572540055904747
$ws.Range("I3").Value = 0.1296533361885828
$ws.Range("J3").Value = 0.1225410148003458
$ws.Range("K3").Value = 0.1225410148003458
$ws.Range("L3").Value = 0.1225410148003458
$ws.Range("M3").Value = 0.1223842413860887
$ws.Range("N3").Value = 0.1223592463389019
$ws.Range("O3").Value = 0.1223592463389019
$ws.Range("P3").Value = 0.1223592463389019
$ws.Range("Q3").Value = 0.1223592463389019
$ws.Range("R3").Value = 0.1223400614984771
$ws.Range("S3").Value = 0.1223400614984771
$ws.Range("T3").Value = 0.1223400614984771
$ws.Range("U3").Value = 0.1223264046107718
$ws.Range("V3").Value = 0.1223264046107718
$ws.Range("W3").Value = 0.1223264046107718
$ws.Range("X3").Value = 0.1223264046107718
$ws.Range("Y3").Value = 0.1223264046107718
$ws.Range("C4").Value = 1.62199854850769
$ws.Range("E4").Value = 5865.680313217402
$ws.Range("F4").Value = 0.2026559610471038
$ws.Range("G4").Value = 0.1808866912957106
$ws.Range("H4").Value = 0.1392480442906786
$ws.Range("I4").Value = 0.133118117096391
$ws.Range("J4").Value = 0.1244335591057
$ws.Range("K4").Value = 0.1225175554845055
$ws.Range("L4").Value = 0.1225175554845055
$ws.Range("M4").Value = 0.1224468796649217
$ws.Range("N4").Value = 0.1224468796649217
$ws.Range("O4").Value = 0.1224468796649217
$ws.Range("P4").Value = 0.1223777202651894
$ws.Range("Q4").Value = 0.1223777202651894
$ws.Range("R4").Value = 0.1223646985842293
$ws.Range("S4").Value = 0.1223646985842293
$ws.Range("T4").Value = 0.1223411227760697
$ws.Range("U4").Value = 0.1223411227760697
$ws.Range("V4").Value = 0.1223407468463431
$ws.Range("W4").Value = 0.1223407468463431
$ws.Range("X4").Value = 0.1223407468463431
$ws.Range("Y4").Value = 0.1223407468463431
$ws.Range("C5").Value = 1.520998239517212
$ws.Range("E5").Value = 5873.581514082798
$ws.Range("F5").Value = 0.2029007438808794
$ws.Range("G5").Value = 0.179815338189779
$ws.Range("H5").Value = 0.1621576929603673
$ws.Range("I5").Value = 0.1441150379024041
$ws.Range("J5").Value = 0.1313006451925719
$ws.Range("K5").Value = 0.1249540660136087
$ws.Range("L5").Value = 0.1231990279500537
$ws.Range("M5").Value = 0.1227704315216329
$ws.Range("N5").Value = 0.1226057157560983
$ws.Range("O5").Value = 0.1225062365900245
$ws.Range("P5").Value = 0.1225062365900245
$ws.Range("Q5").Value = 0.1225062365900245
$ws.Range("R5").Value = 0.1225062365900245
$ws.Range("S5").Value = 0.1225062365900245
$ws.Range("T5").Value = 0.1225023975935622
$ws.Range("U5").Value = 0.1225023975935622
$ws.Range("V5").Value = 0.1225023975935622
$ws.Range("W5").Value = 0.1225023975935622
$ws.Range("X5").Value = 0.1225023975935622
$ws.Range("Y5").Value = 0.1224947663563898
$ws.Range("C6").Value = 1.627990961074829
$ws.Range("E6").Value = 5866.174947578178
$ws.Range("F6").Value = 0.2035703488087754
$ws.Range("G6").Value = 0.1799773616442631
$ws.Range("H6").Value = 0.1510611918736628
$ws.Range("I6").Value = 0.1356651109721328
$ws.Range("J6").Value = 0.1261843268671606
$ws.Range("K6").Value = 0.1237154834177151
$ws.Range("L6").Value = 0.1225052702110248
$ws.Range("M6").Value = 0.1225052702110248
$ws.Range("N6").Value = 0.1223503888416799
$ws.Range("O6").Value = 0.1223503888416799
$ws.Range("P6").Value = 0.1223503888416799
$ws.Range("Q6").Value = 0.1223503888416799
$ws.Range("R6").Value = 0.1223503888416799
$ws.Range("S6").Value = 0.1223503888416799
$ws.Range("T6").Value = 0.1223503888416799
$ws.Range("U6").Value = 0.1223503888416799
$ws.Range("V6").Value = 0.1223503888416799
$ws.Range("W6").Value = 0.1223503888416799
$ws.Range("X6").Value = 0.1223503888416799
$ws.Range("Y6").Value = 0.1223503888416799
$ws.Range("C7").Value = 1.637000322341919
$ws.Range("E7").Value = 5878.233029098237
$ws.Range("F7").Value = 0.2042291535872248
$ws.Range("G7").Value = 0.1761977087256589
$ws.Range("H7").Value = 0.1559592010701578
$ws.Range("I7").Value = 0.1406553088723099
$ws.Range("J7").Value = 0.1271996637517617
$ws.Range("K7").Value = 0.1252265714247468
$ws.Range("L7").Value = 0.1235457448161584
$ws.Range("M7").Value = 0.1232730383909678
$ws.Range("N7").Value = 0.1230048810309338
$ws.Range("O7").Value = 0.1227807866176264
$ws.Range("P7").Value = 0.1227761412609643
$ws.Range("Q7").Value = 0.1227232521035001
$ws.Range("R7").Value = 0.1226335217301097
$ws.Range("S7").Value = 0.1226335217301097
$ws.Range("T7").Value = 0.1226335217301097
$ws.Range("U7").Value = 0.1226141484697698
$ws.Range("V7").Value = 0.1226141484697698
$ws.Range("W7").Value = 0.1226098879784337
$ws.Range("X7").Value = 0.1226098879784337
$ws.Range("Y7").Value = 0.1225854391637083
$ws.Range("C8").Value = 1.554998636245728
$ws.Range("E8").Value = 5865.079115156853
$ws.Range("F8").Value = 0.2064886654668076
$ws.Range("G8").Value = 0.1810726185271584
$ws.Range("H8").Value = 0.1633104977608255
$ws.Range("I8").Value = 0.1446886916624099
$ws.Range("J8").Value = 0.1348851723857546
$ws.Range("K8").Value = 0.1275620079508919
$ws.Range("L8").Value = 0.1229204046266436
$ws.Range("M8").Value = 0.1223612013327923
$ws.Range("N8").Value = 0.1223612013327923
$ws.Range("O8").Value = 0.1223612013327923
$ws.Range("P8").Value = 0.1223290275859035
$ws.Range("Q8").Value = 0.1223290275859035
$ws.Range("R8").Value = 0.1223290275859035
$ws.Range("S8").Value = 0.1223290275859035
$ws.Range("T8").Value = 0.1223290275859035
$ws.Range("U8").Value = 0.1223290275859035
$ws.Range("V8").Value = 0.1223290275859035
$ws.Range("W8").Value = 0.1223290275859035
$ws.Range("X8").Value = 0.1223290275859035
$ws.Range("Y8").Value = 0.1223290275859035
$ws.Range("C9").Value = 1.615995407104492
$ws.Range("E9").Value = 5865.36718856265
$ws.Range("F9").Value = 0.2027056132202467
$ws.Range("G9").Value = 0.1756745569979548
$ws.Range("H9").Value = 0.1415938622175304
$ws.Range("I9").Value = 0.1317587315765017
$ws.Range("J9").Value = 0.1248666469439041
$ws.Range("K9").Value = 0.1225202993754856
$ws.Range("L9").Value = 0.1225202993754856
$ws.Range("M9").Value = 0.1223570921748656
$ws.Range("N9").Value = 0.1223570921748656
$ws.Range("O9").Value = 0.1223570921748656
$ws.Range("P9").Value = 0.1223570921748656
$ws.Range("Q9").Value = 0.1223346430519035
$ws.Range("R9").Value = 0.1223346430519035
$ws.Range("S9").Value = 0.1223346430519035
$ws.Range("T9").Value = 0.1223346430519035
$ws.Range("U9").Value = 0.1223346430519035
$ws.Range("V9").Value = 0.1223346430519035
$ws.Range("W9").Value = 0.1223346430519035
$ws.Range("X9").Value = 0.1223346430519035
$ws.Range("Y9").Value = 0.1223346430519035
$ws.Range("C10").Value = 1.528991222381592
$ws.Range("E10").Value = 5865.049372424783
$ws.Range("F10").Value = 0.2029433039433338
$ws.Range("G10").Value = 0.1829823894811923
$ws.Range("H10").Value = 0.1540499914184416
$ws.Range("I10").Value = 0.1326217124615351
$ws.Range("J10").Value = 0.1228119967476534
$ws.Range("K10").Value = 0.1225781785628779
$ws.Range("L10").Value = 0.1225781785628779
$ws.Range("M10").Value = 0.1225781785628779
$ws.Range("N10").Value = 0.1224894136503342
$ws.Range("O10").Value = 0.1223624131283769
$ws.Range("P10").Value = 0.1223324598295306
$ws.Range("Q10").Value = 0.1223324598295306
$ws.Range("R10").Value = 0.1223324598295306
$ws.Range("S10").Value = 0.1223324598295306
$ws.Range("T10").Value = 0.1223324598295306
$ws.Range("U10").Value = 0.1223324598295306
$ws.Range("V10").Value = 0.1223324598295306
$ws.Range("W10").Value = 0.1223324598295306
$ws.Range("X10").Value = 0.1223324598295306
$ws.Range("Y10").Value = 0.1223284478055513
$ws.Range("C11").Value = 1.622012615203857
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 5864.938817668983
$ws.Range("F11").Value = 0.2044119438235334
$ws.Range("G11").Value = 0.1787928460987157
$ws.Range("H11").Value = 0.1497953860957334
$ws.Range("I11").Value = 0.1363786362240503
$ws.Range("J11").Value = 0.1283876542295382
$ws.Range("K11").Value = 0.1245204671907076
$ws.Range("L11").Value = 0.1223262927420854
$ws.Range("M11").Value = 0.1223262927420854
$ws.Range("N11").Value = 0.1223262927420854
$ws.Range("O11").Value = 0.1223262927420854
$ws.Range("P11").Value = 0.1223262927420854
$ws.Range("Q11").Value = 0.1223262927420854
$ws.Range("R11").Value = 0.1223262927420854
$ws.Range("S11").Value = 0.1223262927420854
$ws.Range("T11").Value = 0.1223262927420854
$ws.Range("U11").Value = 0.1223262927420854
$ws.Range("V11").Value = 0.1223262927420854
$ws.Range("W11").Value = 0.1223262927420854
$ws.Range("X11").Value = 0.1223262927420854
$ws.Range("Y11").Value = 0.1223262927420854
